# This deck's main slide theme (the DrawingML theme bound to the single
# slide master / all slides+layouts) currently uses the "Integral" / "Red
# Violet" colour scheme. The authored change repaints it with the stock
# PowerPoint "Office Theme" colour scheme (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink), which is what the slide master's theme part looked
# like before the swap documented in the commit.
#
# PowerPoint's ColorScheme/RGB values are packed as 0x00BBGGRR (i.e. the
# decimal value of R | (G<<8) | (B<<16)), not 0x00RRGGBB, so each target
# hex colour below is expanded into R/G/B and recombined in that order.

function ToPpRgb([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# Index order matches PowerPoint's ColorSchemeIndex enum:
# 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
# 8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officeThemeColors = @(
    (ToPpRgb 0x00 0x00 0x00),  # dk1      000000
    (ToPpRgb 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (ToPpRgb 0x44 0x54 0x6A),  # dk2      44546A
    (ToPpRgb 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (ToPpRgb 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (ToPpRgb 0xED 0x7D 0x31),  # accent2  ED7D31
    (ToPpRgb 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (ToPpRgb 0xFF 0xC0 0x00),  # accent4  FFC000
    (ToPpRgb 0x44 0x72 0xC4),  # accent5  4472C4
    (ToPpRgb 0x70 0xAD 0x47),  # accent6  70AD47
    (ToPpRgb 0x05 0x63 0xC1),  # hlink    0563C1
    (ToPpRgb 0x95 0x4F 0x72)   # folHlink 954F72
)

# Apply through the slide-level theme colour scheme, which is backed by
# the presentation's single DrawingML theme part shared by every slide,
# the slide master and every layout.
$slide = $p.Slides.Item(1)
$colorScheme = $slide.ThemeColorScheme

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
